$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H113").Value = 144844.42
$ws.Range("I113").Value = 168817.5
$ws.Range("J113").Value = 1006
$ws.Range("K113").Value = 168817.5
$ws.Range("L113").Value = 1006
$ws.Range("M113").Value = -165563.5
$ws.Range("N113").Value = -7514
$ws.Range("H129").Value = 2273.7612
$ws.Range("I129").Value = 5263.7144
$ws.Range("J129").Value = 908.7826
$ws.Range("K129").Value = 15791.1432
$ws.Range("L129").Value = 2726.3478
$ws.Range("M129").Value = -10791.1432
$ws.Range("N129").Value = -12726.3478
$ws.Range("H137").Value = 1589.0588
$ws.Range("I137").Value = 1593.8572
$ws.Range("J137").Value = 1566.6666
$ws.Range("K137").Value = 4781.571599999999
$ws.Range("L137").Value = 4699.9998
$ws.Range("M137").Value = -2231.571599999999
$ws.Range("N137").Value = -9799.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28736.084
$ws.Range("I32").Value = 8567.728999999999
$ws.Range("J32").Value = 120269.38
$ws.Range("K32").Value = 8567.728999999999
$ws.Range("L32").Value = 120269.38
$ws.Range("M32").Value = -8280.728999999999
$ws.Range("N32").Value = -120843.38
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H112").Value = 9347.833000000001
$ws.Range("J112").Value = 9347.833000000001
$ws.Range("L112").Value = 9347.833000000001
$ws.Range("N112").Value = -12301.833
$ws.Range("H132").Value = 3554.4412
$ws.Range("I132").Value = 3094.5862
$ws.Range("J132").Value = 6221.6
$ws.Range("K132").Value = 9283.758600000001
$ws.Range("L132").Value = 18664.8
$ws.Range("M132").Value = -6753.758600000001
$ws.Range("N132").Value = -23724.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 26600
$ws.Range("J68").Value = 26600
$ws.Range("L68").Value = 26600
$ws.Range("N68").Value = -28222
$ws.Range("H71").Value = 26600
$ws.Range("J71").Value = 26600
$ws.Range("L71").Value = 79800
$ws.Range("N71").Value = -87912

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15249.25
$ws.Range("J99").Value = 16859.143
$ws.Range("L99").Value = 16859.143
$ws.Range("N99").Value = -19855.143
$ws.Range("H126").Value = 15249.25
$ws.Range("J126").Value = 16859.143
$ws.Range("L126").Value = 50577.429
$ws.Range("N126").Value = -55517.429
$ws.Range("H134").Value = 1656.6875
$ws.Range("I134").Value = 1573.909
$ws.Range("J134").Value = 1838.8
$ws.Range("K134").Value = 4721.727000000001
$ws.Range("L134").Value = 5516.4
$ws.Range("M134").Value = -2186.727000000001
$ws.Range("N134").Value = -10586.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 812
$ws.Range("J20").Value = 60
$ws.Range("L20").Value = 180
$ws.Range("N20").Value = -634
$ws.Range("H22").Value = 7537.125
$ws.Range("J22").Value = 9699.5
$ws.Range("L22").Value = 29098.5
$ws.Range("N22").Value = -29436.5
$ws.Range("H24").Value = 630
$ws.Range("J24").Value = 630
$ws.Range("L24").Value = 1890
$ws.Range("N24").Value = -2350
$ws.Range("H25").Value = 600
$ws.Range("I25").Value = 450
$ws.Range("J25").Value = 675
$ws.Range("K25").Value = 1350
$ws.Range("L25").Value = 2025
$ws.Range("M25").Value = -1181
$ws.Range("N25").Value = -2363
$ws.Range("H27").Value = 7537.125
$ws.Range("J27").Value = 9699.5
$ws.Range("L27").Value = 29098.5
$ws.Range("N27").Value = -29302.5
$ws.Range("H29").Value = 174.78572
$ws.Range("I29").Value = 46.75
$ws.Range("K29").Value = 140.25
$ws.Range("M29").Value = 136.75
$ws.Range("H30").Value = 600
$ws.Range("I30").Value = 450
$ws.Range("J30").Value = 675
$ws.Range("K30").Value = 1350
$ws.Range("L30").Value = 2025
$ws.Range("M30").Value = -1248
$ws.Range("N30").Value = -2229
$ws.Range("H34").Value = 938.2941
$ws.Range("J34").Value = 1133.0834
$ws.Range("L34").Value = 3399.2502
$ws.Range("N34").Value = -3567.2502
$ws.Range("H39").Value = 400
$ws.Range("J39").Value = 300
$ws.Range("L39").Value = 900
$ws.Range("N39").Value = -1488
$ws.Range("H44").Value = 696.5
$ws.Range("I44").Value = 198.66667
$ws.Range("J44").Value = 2190
$ws.Range("K44").Value = 596.00001
$ws.Range("L44").Value = 6570
$ws.Range("M44").Value = -198.00001
$ws.Range("N44").Value = -7366
$ws.Range("H46").Value = 126876
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 169001.33
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 507003.99
$ws.Range("M46").Value = -1409
$ws.Range("N46").Value = -507185.99
$ws.Range("H55").Value = 17083.334
$ws.Range("J55").Value = 12837.5
$ws.Range("L55").Value = 38512.5
$ws.Range("N55").Value = -38866.5
$ws.Range("H57").Value = 1000
$ws.Range("I57").Value = 1000
$ws.Range("K57").Value = 3000
$ws.Range("M57").Value = -2441
$ws.Range("H113").Value = 779.34283
$ws.Range("I113").Value = 1380.5454
$ws.Range("J113").Value = 503.79166
$ws.Range("K113").Value = 4141.6362
$ws.Range("L113").Value = 1511.37498
$ws.Range("M113").Value = -1971.6362
$ws.Range("N113").Value = -5851.374980000001
$ws.Range("H131").Value = 1058.1111
$ws.Range("J131").Value = 1063.1549
$ws.Range("L131").Value = 3189.4647
$ws.Range("N131").Value = -13269.4647
$ws.Range("H132").Value = 1707.0555
$ws.Range("I132").Value = 900.4737
$ws.Range("J132").Value = 2608.5293
$ws.Range("K132").Value = 8104.263300000001
$ws.Range("L132").Value = 23476.7637
$ws.Range("M132").Value = -5574.263300000001
$ws.Range("N132").Value = -28536.7637

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H126").Value = 2587.4285
$ws.Range("I126").Value = 2684
$ws.Range("J126").Value = 2233.3333
$ws.Range("K126").Value = 8052
$ws.Range("L126").Value = 6699.999899999999
$ws.Range("M126").Value = -5582
$ws.Range("N126").Value = -11639.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 42678.72
$ws.Range("I40").Value = 69624.53
$ws.Range("J40").Value = 2260
$ws.Range("K40").Value = 69624.53
$ws.Range("L40").Value = 2260
$ws.Range("M40").Value = -69488.53
$ws.Range("N40").Value = -2532
$ws.Range("H139").Value = 48976.668
$ws.Range("I139").Value = 37000
$ws.Range("J139").Value = 54965
$ws.Range("K139").Value = 37000
$ws.Range("L139").Value = 54965
$ws.Range("M139").Value = -31860
$ws.Range("N139").Value = -65245

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 0
$ws.Range("N126").ClearContents()
